$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-03 15:41:15"

# Sheets 2 ("Главные") and 3 ("Линейные") both carry an as_of_utc
# timestamp column (AA) for rows 2-26 that needs to be refreshed.
$sheetIndexes = @(2, 3)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column 27 = AA
        $cell.Value = $newTimestamp
    }
}
